$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "state" column header to "country"
$ws.Range("A1").Value = "country"

# Move the active selection from E2 to A2
$ws.Range("A2").Select()
